# Trade #39 closed at 2026-02-17 15:23:44 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.13   # Total P&L %
$summary.Range("B6").Value = 39      # Total Trades
$summary.Range("B9").Value = 28.21   # Win Rate %

# --- Strategy Status sheet (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 39       # Trades
$status.Range("G4").Value = 28.21    # Win Rate %

# --- New trade row (#39), appended to both "All Trades" and "MarketMaking" ---
foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A40").Value = 39

    # Force the date-looking / time-looking strings to stay plain text
    # (matching the rest of the column) instead of being auto-converted
    # to a date/time serial value.
    $ws.Range("B40").NumberFormat = "@"
    $ws.Range("B40").Value = "2026-02-17"
    $ws.Range("B40").Style = "Normal"

    $ws.Range("C40").Value = "15:23:38"
    $ws.Range("D40").Value = "MarketMaking"
    $ws.Range("E40").Value = "DOWN"
    $ws.Range("F40").Value = 0.59
    $ws.Range("G40").Value = 0.59
    $ws.Range("H40").Value = "CLOSED"
    $ws.Range("I40").Value = 0
    $ws.Range("J40").Value = 0
    $ws.Range("K40").Value = 99.73999999999999
    $ws.Range("L40").Value = 0
    $ws.Range("M40").Value = 0
    $ws.Range("N40").Value = 0.6
    $ws.Range("O40").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P40").Value = "early_exit"
    $ws.Range("Q40").Value = 0.14
}
